$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.080435
$ws.Range("H2").Value = 24.241305
$ws.Range("I2").Value = 0.1496988574979475
$ws.Range("J2").Value = 0.1496988574979476
$ws.Range("M2").Value = 2.590549
$ws.Range("N2").Value = 7.771647
$ws.Range("O2").Value = 0.08453422544559429
$ws.Range("P2").Value = 0.0845342254455943
$ws.Range("Q2").Value = 20.932762808815
$ws.Range("R2").Value = 188.394865279335
$ws.Range("S2").Value = 0.01265467696867939
$ws.Range("T2").Value = 0.01265467696867939
$ws.Range("G3").Value = 8.080435
$ws.Range("H3").Value = 24.241305
$ws.Range("I3").Value = 0.1496988574979475
$ws.Range("J3").Value = 0.1496988574979476
$ws.Range("O3").Value = 0.1976233469596758
$ws.Range("P3").Value = 0.1976233469596758
$ws.Range("Q3").Value = 48.93642339047
$ws.Range("R3").Value = 440.42781051423
$ws.Range("S3").Value = 0.02958398925478396
$ws.Range("T3").Value = 0.02958398925478396
$ws.Range("G4").Value = 8.080435
$ws.Range("H4").Value = 24.241305
$ws.Range("I4").Value = 0.1496988574979475
$ws.Range("J4").Value = 0.1496988574979476
$ws.Range("M4").Value = 12.64302866666667
$ws.Range("N4").Value = 37.929086
$ws.Range("O4").Value = 0.4125645319286033
$ws.Range("P4").Value = 0.4125645319286034
$ws.Range("Q4").Value = 102.1611713441367
$ws.Range("R4").Value = 919.45054209723
$ws.Range("S4").Value = 0.06176043907388742
$ws.Range("T4").Value = 0.06176043907388743
$ws.Range("G5").Value = 8.080435
$ws.Range("H5").Value = 24.241305
$ws.Range("I5").Value = 0.1496988574979475
$ws.Range("J5").Value = 0.1496988574979476
$ws.Range("M5").Value = 3.370263
$ws.Range("N5").Value = 10.110789
$ws.Range("O5").Value = 0.1099776812764186
$ws.Range("P5").Value = 0.1099776812764186
$ws.Range("Q5").Value = 27.233191104405
$ws.Range("R5").Value = 245.098719939645
$ws.Range("S5").Value = 0.01646353323735328
$ws.Range("T5").Value = 0.01646353323735329
$ws.Range("G6").Value = 8.080435
$ws.Range("H6").Value = 24.241305
$ws.Range("I6").Value = 0.1496988574979475
$ws.Range("J6").Value = 0.1496988574979476
$ws.Range("M6").Value = 5.984969666666667
$ws.Range("N6").Value = 17.954909
$ws.Range("O6").Value = 0.1953002143897079
$ws.Range("P6").Value = 0.1953002143897079
$ws.Range("Q6").Value = 48.36115836847167
$ws.Range("R6").Value = 435.250425316245
$ws.Range("S6").Value = 0.02923621896324348
$ws.Range("T6").Value = 0.02923621896324349
$ws.Range("I7").Value = 0.2404784903431001
$ws.Range("J7").Value = 0.2404784903431001
$ws.Range("M7").Value = 2.590549
$ws.Range("N7").Value = 7.771647
$ws.Range("O7").Value = 0.08453422544559429
$ws.Range("P7").Value = 0.0845342254455943
$ws.Range("Q7").Value = 33.626704192068
$ws.Range("R7").Value = 302.640337728612
$ws.Range("S7").Value = 0.0203286629174798
$ws.Range("T7").Value = 0.0203286629174798
$ws.Range("I8").Value = 0.2404784903431001
$ws.Range("J8").Value = 0.2404784903431001
$ws.Range("O8").Value = 0.1976233469596758
$ws.Range("P8").Value = 0.1976233469596758
$ws.Range("S8").Value = 0.04752416413341354
$ws.Range("T8").Value = 0.04752416413341354
$ws.Range("I9").Value = 0.2404784903431001
$ws.Range("J9").Value = 0.2404784903431001
$ws.Range("M9").Value = 12.64302866666667
$ws.Range("N9").Value = 37.929086
$ws.Range("O9").Value = 0.4125645319286033
$ws.Range("P9").Value = 0.4125645319286034
$ws.Range("Q9").Value = 164.113238184584
$ws.Range("R9").Value = 1477.019143661256
$ws.Range("S9").Value = 0.09921289580729827
$ws.Range("T9").Value = 0.09921289580729828
$ws.Range("I10").Value = 0.2404784903431001
$ws.Range("J10").Value = 0.2404784903431001
$ws.Range("M10").Value = 3.370263
$ws.Range("N10").Value = 10.110789
$ws.Range("O10").Value = 0.1099776812764186
$ws.Range("P10").Value = 0.1099776812764186
$ws.Range("Q10").Value = 43.74780671991601
$ws.Range("R10").Value = 393.7302604792441
$ws.Range("S10").Value = 0.02644726676478778
$ws.Range("T10").Value = 0.02644726676478778
$ws.Range("I11").Value = 0.2404784903431001
$ws.Range("J11").Value = 0.2404784903431001
$ws.Range("M11").Value = 5.984969666666667
$ws.Range("N11").Value = 17.954909
$ws.Range("O11").Value = 0.1953002143897079
$ws.Range("P11").Value = 0.1953002143897079
$ws.Range("Q11").Value = 77.68809027719603
$ws.Range("R11").Value = 699.1928124947641
$ws.Range("S11").Value = 0.04696550072012075
$ws.Range("T11").Value = 0.04696550072012076
$ws.Range("G12").Value = 15.25749233333333
$ws.Range("H12").Value = 45.772477
$ws.Range("I12").Value = 0.2826616599952471
$ws.Range("J12").Value = 0.2826616599952471
$ws.Range("M12").Value = 2.590549
$ws.Range("N12").Value = 7.771647
$ws.Range("O12").Value = 0.08453422544559429
$ws.Range("P12").Value = 0.0845342254455943
$ws.Range("Q12").Value = 39.52528150662433
$ws.Range("R12").Value = 355.727533559619
$ws.Range("S12").Value = 0.02389458449086413
$ws.Range("T12").Value = 0.02389458449086414
$ws.Range("G13").Value = 15.25749233333333
$ws.Range("H13").Value = 45.772477
$ws.Range("I13").Value = 0.2826616599952471
$ws.Range("J13").Value = 0.2826616599952471
$ws.Range("O13").Value = 0.1976233469596758
$ws.Range("P13").Value = 0.1976233469596758
$ws.Range("Q13").Value = 92.40184528442468
$ws.Range("R13").Value = 831.6166075598221
$ws.Range("S13").Value = 0.05586054330543864
$ws.Range("T13").Value = 0.05586054330543864
$ws.Range("G14").Value = 15.25749233333333
$ws.Range("H14").Value = 45.772477
$ws.Range("I14").Value = 0.2826616599952471
$ws.Range("J14").Value = 0.2826616599952471
$ws.Range("M14").Value = 12.64302866666667
$ws.Range("N14").Value = 37.929086
$ws.Range("O14").Value = 0.4125645319286033
$ws.Range("P14").Value = 0.4125645319286034
$ws.Range("Q14").Value = 192.9009129517802
$ws.Range("R14").Value = 1736.108216566022
$ws.Range("S14").Value = 0.1166161754501011
$ws.Range("T14").Value = 0.1166161754501011
$ws.Range("G15").Value = 15.25749233333333
$ws.Range("H15").Value = 45.772477
$ws.Range("I15").Value = 0.2826616599952471
$ws.Range("J15").Value = 0.2826616599952471
$ws.Range("M15").Value = 3.370263
$ws.Range("N15").Value = 10.110789
$ws.Range("O15").Value = 0.1099776812764186
$ws.Range("P15").Value = 0.1099776812764186
$ws.Range("Q15").Value = 51.421761883817
$ws.Range("R15").Value = 462.795856954353
$ws.Range("S15").Value = 0.03108647395202069
$ws.Range("T15").Value = 0.03108647395202069
$ws.Range("G16").Value = 15.25749233333333
$ws.Range("H16").Value = 45.772477
$ws.Range("I16").Value = 0.2826616599952471
$ws.Range("J16").Value = 0.2826616599952471
$ws.Range("M16").Value = 5.984969666666667
$ws.Range("N16").Value = 17.954909
$ws.Range("O16").Value = 0.1953002143897079
$ws.Range("P16").Value = 0.1953002143897079
$ws.Range("Q16").Value = 91.31562880439922
$ws.Range("R16").Value = 821.840659239593
$ws.Range("S16").Value = 0.05520388279682246
$ws.Range("T16").Value = 0.05520388279682247
$ws.Range("G17").Value = 4.142925
$ws.Range("H17").Value = 12.428775
$ws.Range("I17").Value = 0.07675219702895753
$ws.Range("J17").Value = 0.07675219702895753
$ws.Range("M17").Value = 2.590549
$ws.Range("N17").Value = 7.771647
$ws.Range("O17").Value = 0.08453422544559429
$ws.Range("P17").Value = 0.0845342254455943
$ws.Range("Q17").Value = 10.732450215825
$ws.Range("R17").Value = 96.592051942425
$ws.Range("S17").Value = 0.006488187527090568
$ws.Range("T17").Value = 0.006488187527090569
$ws.Range("G18").Value = 4.142925
$ws.Range("H18").Value = 12.428775
$ws.Range("I18").Value = 0.07675219702895753
$ws.Range("J18").Value = 0.07675219702895753
$ws.Range("O18").Value = 0.1976233469596758
$ws.Range("P18").Value = 0.1976233469596758
$ws.Range("Q18").Value = 25.09022495385
$ws.Range("R18").Value = 225.81202458465
$ws.Range("S18").Value = 0.01516802606337108
$ws.Range("T18").Value = 0.01516802606337108
$ws.Range("G19").Value = 4.142925
$ws.Range("H19").Value = 12.428775
$ws.Range("I19").Value = 0.07675219702895753
$ws.Range("J19").Value = 0.07675219702895753
$ws.Range("M19").Value = 12.64302866666667
$ws.Range("N19").Value = 37.929086
$ws.Range("O19").Value = 0.4125645319286033
$ws.Range("P19").Value = 0.4125645319286034
$ws.Range("Q19").Value = 52.37911953885
$ws.Range("R19").Value = 471.41207584965
$ws.Range("S19").Value = 0.0316652342417438
$ws.Range("T19").Value = 0.03166523424174381
$ws.Range("G20").Value = 4.142925
$ws.Range("H20").Value = 12.428775
$ws.Range("I20").Value = 0.07675219702895753
$ws.Range("J20").Value = 0.07675219702895753
$ws.Range("M20").Value = 3.370263
$ws.Range("N20").Value = 10.110789
$ws.Range("O20").Value = 0.1099776812764186
$ws.Range("P20").Value = 0.1099776812764186
$ws.Range("Q20").Value = 13.962746839275
$ws.Range("R20").Value = 125.664721553475
$ws.Range("S20").Value = 0.008441028662115574
$ws.Range("T20").Value = 0.008441028662115576
$ws.Range("G21").Value = 4.142925
$ws.Range("H21").Value = 12.428775
$ws.Range("I21").Value = 0.07675219702895753
$ws.Range("J21").Value = 0.07675219702895753
$ws.Range("M21").Value = 5.984969666666667
$ws.Range("N21").Value = 17.954909
$ws.Range("O21").Value = 0.1953002143897079
$ws.Range("P21").Value = 0.1953002143897079
$ws.Range("Q21").Value = 24.795280456275
$ws.Range("R21").Value = 223.157524106475
$ws.Range("S21").Value = 0.0149897205346365
$ws.Range("T21").Value = 0.01498972053463651
$ws.Range("G22").Value = 13.51654933333334
$ws.Range("H22").Value = 40.549648
$ws.Range("I22").Value = 0.2504087951347477
$ws.Range("J22").Value = 0.2504087951347477
$ws.Range("M22").Value = 2.590549
$ws.Range("N22").Value = 7.771647
$ws.Range("O22").Value = 0.08453422544559429
$ws.Range("P22").Value = 0.0845342254455943
$ws.Range("Q22").Value = 35.01528335891734
$ws.Range("R22").Value = 315.137550230256
$ws.Range("S22").Value = 0.02116811354148039
$ws.Range("T22").Value = 0.0211681135414804
$ws.Range("G23").Value = 13.51654933333334
$ws.Range("H23").Value = 40.549648
$ws.Range("I23").Value = 0.2504087951347477
$ws.Range("J23").Value = 0.2504087951347477
$ws.Range("O23").Value = 0.1976233469596758
$ws.Range("P23").Value = 0.1976233469596758
$ws.Range("Q23").Value = 81.85841244365869
$ws.Range("R23").Value = 736.7257119929282
$ws.Range("S23").Value = 0.04948662420266862
$ws.Range("T23").Value = 0.04948662420266862
$ws.Range("G24").Value = 13.51654933333334
$ws.Range("H24").Value = 40.549648
$ws.Range("I24").Value = 0.2504087951347477
$ws.Range("J24").Value = 0.2504087951347477
$ws.Range("M24").Value = 12.64302866666667
$ws.Range("N24").Value = 37.929086
$ws.Range("O24").Value = 0.4125645319286033
$ws.Range("P24").Value = 0.4125645319286034
$ws.Range("Q24").Value = 170.8901206957476
$ws.Range("R24").Value = 1538.011086261728
$ws.Range("S24").Value = 0.1033097873555727
$ws.Range("T24").Value = 0.1033097873555727
$ws.Range("G25").Value = 13.51654933333334
$ws.Range("H25").Value = 40.549648
$ws.Range("I25").Value = 0.2504087951347477
$ws.Range("J25").Value = 0.2504087951347477
$ws.Range("M25").Value = 3.370263
$ws.Range("N25").Value = 10.110789
$ws.Range("O25").Value = 0.1099776812764186
$ws.Range("P25").Value = 0.1099776812764186
$ws.Range("Q25").Value = 45.55432610580801
$ws.Range("R25").Value = 409.9889349522721
$ws.Range("S25").Value = 0.02753937866014129
$ws.Range("T25").Value = 0.02753937866014129
$ws.Range("G26").Value = 13.51654933333334
$ws.Range("H26").Value = 40.549648
$ws.Range("I26").Value = 0.2504087951347477
$ws.Range("J26").Value = 0.2504087951347477
$ws.Range("M26").Value = 5.984969666666667
$ws.Range("N26").Value = 17.954909
$ws.Range("O26").Value = 0.1953002143897079
$ws.Range("P26").Value = 0.1953002143897079
$ws.Range("Q26").Value = 80.89613775800358
$ws.Range("R26").Value = 728.0652398220321
$ws.Range("S26").Value = 0.04890489137488466
$ws.Range("T26").Value = 0.04890489137488466
